# Generate Report for Handback
#
# The handback transform failed for the 7c831e84-... source file in both the
# zh-cn and de-de locale sheets. This updates:
#   - the Status cell (col C, row 3) on the zh-cn / de-de sheets (and the
#     mirrored Overview columns that show the same status) from
#     "Ready for handoff" to "Handback transform failed"
#   - the Error Detail cell (col P, row 3) on the zh-cn / de-de sheets with
#     the transform-failure message
#   - widens the Error Detail column so the longer message is readable

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Status column (every cell that previously read "Ready for handoff")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Error Detail column for the failed handback, per locale
$zhcn.Range("P3").Value = "Handback file name: wsnwuf2w.v3d is different with handoff file name: 7c831e84-5a9f-4373-99f5-83a5673f4749.b63ffbb5156f5652266d7444e4cf6b4487c16f3e.zh-cn."
$dede.Range("P3").Value = "Handback file name: wsnwuf2w.v3d is different with handoff file name: 7c831e84-5a9f-4373-99f5-83a5673f4749.b63ffbb5156f5652266d7444e4cf6b4487c16f3e.de-de."

# Widen the Error Detail column (P) on both locale sheets to fit the new
# message (stored column width ends up 40 characters wide).
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
